$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Munka1")

# New timing data (sparse model changed -> re-ran testing) for Sheet1 B2:H41
$data = @(
    @(0.38100000000000001,2.9249999999999998,0.59,1.4990000000000001,5.3940000000000001,0.89,39),
    @(0.41299999999999998,2.4350000000000001,0.59699999999999998,1.306,4.7519999999999998,0.90500000000000003,27),
    @(0.32200000000000001,2.363,0.61899999999999999,1.2629999999999999,4.5670000000000002,0.70299999999999996,32),
    @(0.55800000000000005,2.23,0.56999999999999995,1.5860000000000001,4.9450000000000003,0.86299999999999999,33),
    @(0.32500000000000001,2.2829999999999999,0.57099999999999995,1.319,4.4980000000000002,1,34),
    @(0.38400000000000001,2.2080000000000002,0.61599999999999999,3.677,6.8849999999999998,0.72499999999999998,35),
    @(0.29299999999999998,2.3679999999999999,0.59599999999999997,1.4690000000000001,4.726,0.74199999999999999,36),
    @(0.311,3.2229999999999999,0.56100000000000005,1.3919999999999999,5.4859999999999998,0.78900000000000003,31),
    @(0.30299999999999999,2.3180000000000001,0.59099999999999997,1.5609999999999999,4.7729999999999997,1,40),
    @(0.58799999999999997,2.2519999999999998,0.64300000000000002,1.242,4.7249999999999996,0.87,22),
    @(1,2.3279999999999998,0.53900000000000003,1.587,5.4539999999999997,0.86499999999999999,40),
    @(0.29499999999999998,2.5139999999999998,0.58899999999999997,1.5149999999999999,4.9139999999999997,0.78800000000000003,34),
    @(0.35099999999999998,2.7450000000000001,0.58799999999999997,1.425,5.109,1,28),
    @(0.41199999999999998,2.9809999999999999,0.57499999999999996,1.798,5.7670000000000003,0.68799999999999994,28),
    @(0.33300000000000002,2.577,0.58399999999999996,1.704,5.1989999999999998,1,28),
    @(0.309,2.5960000000000001,0.58799999999999997,1.2789999999999999,4.7720000000000002,0.93,27),
    @(0.30399999999999999,2.2770000000000001,0.57299999999999995,1.877,5.0309999999999997,0.74099999999999999,28),
    @(0.32100000000000001,2.2959999999999998,0.622,1.552,4.7910000000000004,0.70199999999999996,25),
    @(0.41899999999999998,2.516,0.623,1.2470000000000001,4.8040000000000003,0.7,17),
    @(0.28299999999999997,2.42,0.60599999999999998,1.2290000000000001,4.5380000000000003,0.86399999999999999,29),
    @(0.34499999999999997,3.1949999999999998,0.52900000000000003,1.222,5.2910000000000004,0.94099999999999995,36),
    @(0.36399999999999999,2.177,0.57699999999999996,1.6479999999999999,4.7649999999999997,0.68600000000000005,20),
    @(0.33100000000000002,2.367,0.63800000000000001,2.508,5.8449999999999998,0.75800000000000001,28),
    @(0.308,3.6320000000000001,0.62,1.7589999999999999,6.319,0.9,32),
    @(0.40500000000000003,2.4470000000000001,0.57799999999999996,1.659,5.09,0.63800000000000001,31),
    @(0.29899999999999999,2.2440000000000002,0.51500000000000001,1.1990000000000001,4.2569999999999997,0.73499999999999999,8),
    @(0.45800000000000002,2.4129999999999998,0.56599999999999995,1.798,5.2350000000000003,0.64900000000000002,32),
    @(0.27600000000000002,2.1629999999999998,0.55000000000000004,1.2869999999999999,4.2759999999999998,0.98099999999999998,9),
    @(0.29399999999999998,2.38,0.56000000000000005,1.4350000000000001,4.6680000000000001,0.755,8),
    @(0.36,2.399,0.58199999999999996,1.498,4.8380000000000001,0.83899999999999997,27),
    @(0.35399999999999998,2.266,0.60599999999999998,1.6890000000000001,4.915,0.96,19),
    @(0.314,2.528,0.59799999999999998,1.677,5.117,0.88300000000000001,34),
    @(0.34799999999999998,2.4750000000000001,0.58199999999999996,1.839,5.2430000000000003,0.89900000000000002,42),
    @(0.34200000000000003,2.2349999999999999,0.58499999999999996,1.6579999999999999,4.82,0.65,33),
    @(0.30499999999999999,2.3889999999999998,0.56000000000000005,1.611,4.8650000000000002,0.89500000000000002,39),
    @(0.311,2.754,0.6,1.46,5.1260000000000003,1,31),
    @(0.27300000000000002,2.2770000000000001,0.63100000000000001,1.948,5.1280000000000001,0.76400000000000001,36),
    @(0.28799999999999998,2.371,0.58499999999999996,1.151,4.3949999999999996,1,29),
    @(0.32,2.383,0.54500000000000004,1.3240000000000001,4.5720000000000001,0.90300000000000002,29),
    @(0.46600000000000003,2.2799999999999998,0.58599999999999997,1.4490000000000001,4.7809999999999997,0.88400000000000001,34)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws1.Cells.Item($i + 2, $j + 2).Value = $row[$j]
    }
}

# Sheet2 (Munka1) G5/G7/G9/G11/G13/G15/G17 hold =ROUND(AVERAGE(...),2) formulas
# that recalc automatically once Sheet1 data changes above.

# Reset view / selection state to match the edited workbook
$ws1.Range("A1").Select()
$ws2.Range("G18").Select()
